$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set E column ("y (t/cu.m.)") value = 3.8 for rows 5-25 (row 4 already has it)
for ($r = 5; $r -le 25; $r++) {
    $ws.Cells.Item($r, 5).Value = 3.8
}

# Remove K column (Soil Class) values for rows 15-25 that referenced the
# now-removed "Granite" family soil-class strings
for ($r = 15; $r -le 25; $r++) {
    $ws.Cells.Item($r, 11).Clear()
}

# Update the saved selection to L17
$ws.Range("L17").Select()
